# Generate Report for Handback
# Adds a new handback-status row for file
# f2fa4ce0-bf40-4efe-91dc-85b8937f538d.md (hash 55bf4393a37a88e53ecf3ac86c08ef0c4cfb56cf)
# to every sheet: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$baseName   = "f2fa4ce0-bf40-4efe-91dc-85b8937f538d"
$mdFile     = "$baseName.md"
$hash       = "55bf4393a37a88e53ecf3ac86c08ef0c4cfb56cf"
$zhFile     = "$baseName.$hash.zh-cn.xlf"
$deFile     = "$baseName.$hash.de-de.xlf"

$status     = "Handed back: in sync with en-US"
$reason     = "Include"

$zhHandoffDt  = "2016-02-22 09:10:59"
$zhHandbackDt = "2016-02-22 09:11:42"
$deHandoffDt  = "2016-02-22 09:11:13"
$deHandbackDt = "2016-02-22 09:12:05"

# --------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# --------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$rowOv = 4

$wsOverview.Cells.Item($rowOv, 2).Value2 = $status
$wsOverview.Cells.Item($rowOv, 3).Value2 = $status

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($rowOv, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/$hash/e2e/$mdFile",
    "",
    "",
    $mdFile
) | Out-Null

# --------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Correspond Handoff File |
#   Correspond Handoff Datetime | Target File | Correspond Handback File |
#   Correspond Handback DateTime | Handoff Reason | Dependency From
# --------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$rowZh = 4

$wsZh.Cells.Item($rowZh, 2).Value2 = $status
$wsZh.Cells.Item($rowZh, 4).Value2 = $zhHandoffDt
$wsZh.Cells.Item($rowZh, 7).Value2 = $zhHandbackDt
$wsZh.Cells.Item($rowZh, 8).Value2 = $reason

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item($rowZh, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/$hash/e2e/$mdFile",
    "",
    "",
    $mdFile
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item($rowZh, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/$zhFile",
    "",
    "",
    $zhFile
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item($rowZh, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$hash/e2e/$mdFile",
    "",
    "",
    $mdFile
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Cells.Item($rowZh, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$hash/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/$zhFile",
    "",
    "",
    $zhFile
) | Out-Null

# --------------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn
# --------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$rowDe = 4

$wsDe.Cells.Item($rowDe, 2).Value2 = $status
$wsDe.Cells.Item($rowDe, 4).Value2 = $deHandoffDt
$wsDe.Cells.Item($rowDe, 7).Value2 = $deHandbackDt
$wsDe.Cells.Item($rowDe, 8).Value2 = $reason

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item($rowDe, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/$hash/e2e/$mdFile",
    "",
    "",
    $mdFile
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item($rowDe, 3),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/$deFile",
    "",
    "",
    $deFile
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item($rowDe, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$hash/e2e/$mdFile",
    "",
    "",
    $mdFile
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Cells.Item($rowDe, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$hash/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/$deFile",
    "",
    "",
    $deFile
) | Out-Null
